$d = $word.ActiveDocument

# --- Paragraph 1: title/date line (two runs split by a line break) ---
$d.Content.Find.Execute(
    "🚀המאמר היומי של מייק -11.10.24: ⚡️🚀", $true, $false, $false, $false, $false,
    $true, 1, $false, "🚀המאמר היומי של מייק -10.10.24: ⚡️🚀", 2) | Out-Null

$d.Content.Find.Execute(
    "SELECTIVE ATTENTION IMPROVES TRANSFORMER", $true, $false, $false, $false, $false,
    $true, 1, $false, "DIFFERENTIAL TRANSFORMER", 2) | Out-Null

# --- Paragraph 2 ---
$d.Content.Find.Execute(
    "היום נסקור מאמר המציג רעיון לשיפור הליבה של הטרנספורמים, כלומר מנגנון ה-attention. להבדיל מהמאמר של סקרתי(Selective Transformer) הרעיון כאן די ברור לי מתמטית ולא ולא זיהיתי בו נוסחאות מתמטיות ״מפתיעות״. המאמר של היום מציע שיטה לשיפור ביצועים של הטרנספורמרים ועל הדרך מצליח להקטין את גודל הזכרון הנדרש עבורו. ",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "המאמר הזה עשה הרבה גלים ביומיים האחרונים וזו הסיבה שבחרתי אותו לסקירה היומית שלי. המאמר החזיר אותי 3-4 שנים אחורה לתקופה שבה על בסיס ימי יצאו מאמרים המציעים שכלולים שונים לליבה של הטרנספורמרים כה אהובים עלינו. כמובן אני מתכוון למנגנון ה-attention שמאפשר לנו לכמת קשרים בין הטוקנים השונים בטקסט. ",
    2) | Out-Null

# --- Paragraph 3 ---
$d.Content.Find.Execute(
    "המחברים טוענים (ובצדק) שלפעמים יש טוקנים שלא כדאי לטרוח ולחשב מקדמי attention עבור זוגות מסוימים של הטוקנים. בנוסף ניתן לדעת את זה על ידי הסתכלות על טוקנים ביניהם ואלו באים לפניהם (ההקשר). ֿ",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "המחברים הציעו להחליף את חישוב הסופטמקס הרגיל שיש לנו בטרנספורמרים בהפרש משוקלל (רק הסופטמקס השני משוקלל) של הסופטמקסים. כל סופטמקס מחושב עם מטריצת Q ו-K משלה כאשר המשקול λ של הסופטמקס השני מחושב באופן הבא:  λ = exp(λ_q1 · λ_k1 ) − exp(λ_q2 · λ_k2 ) + λ_init כאשר ",
    2) | Out-Null

# --- Paragraph 4 ---
$d.Content.Find.Execute(
    "המחברים נותנים את הדוגמא הבאה הממחישה את התופעה הזו. נניח שהטוקנים א, ב, ג הוזנו לטרנספורמר. בשכבה כלשהו עם מ attention סטנדרטי, טוקן ב מחליט ״כמה הוא מעוניין לקחת״ מטוקן א (מקדם attention), וטוקן ג יכול להחליט כמה לקרוא מטוקן א, אבל טוקן ב אינו יכול להשפיע על כמה טוקן ג ״לוקח״ מטוקן א. אם טוקן ב קבע שטוקן א אינו רלוונטי או אפילו מטעה לטוקנים עתידיים כמו ג, אין שום דבר שהוא יכול לעשות בשכבה הנתונה כדי לתקן זאת. השיטה המוצעת על ידי המחברים באה לתקן (להקל) את הבעיה הזו.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "λ_q1 , λ_k1 , λ_q2 , λ_k2 ∈  R^d הינם נלמדים ו- ((λ_init = 0.8 − 0.6 × exp(−0.3 · (l − 1, כאשר l זה מספר השכבה (של בלוק הטרנספורמר). אם הנוסחה עבור λ איכשהו מובנת ודי סטנדרטית הנוסחה עבור λ_init נותרת בגדר תעלומה (אלא אם כן זה ניסוי ותהיה רגרסיה של הערכים שהתקבלו עם פונקציה מצורה מסוימת).",
    2) | Out-Null

# --- Paragraph 5 ---
$d.Content.Find.Execute(
    "הרעיון המוצע הוא מאוד אינטואיטיבי ואלגנטי. המחברים מציע להחסיר מווקטור ה-attention (לפני חישוב הסופטמקס) של כל טוקן מטריצת מיסוך נלמדת F. איבר ij במטריצת F (עבור זוג טוקנים i- ו- j) מבטא עד כמה אנו רוצים להקטין את ה-attention בין טוקנים אלו. ערך גבוה של F_ij מסמן לנו שהמודל ״מאמין״ שצריך ״להתעלם מהקשר בין טוקן i ל-j כלומר (אם i > j; מטריצה F הינה מטריצה קוזלית כלומר F_ij =0 אם i<j) אפשר לא לדלג על חישוב מקדם ה-attention ביניהם. ",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "המאמר טוען לשיפור תוצאת אבל הבדיקות נעשו בעיקר למודלים עם 3B פרמטרים. יש גם טענות לקנסול של רעש כלשהו שאני לא בטוח שאני מבין. בקיצר אני קצת סקפטי, מודה….",
    2) | Out-Null

# --- Paragraph 10 (the arxiv link) gets its text swapped in place ---
$d.Content.Find.Execute(
    "https://arxiv.org/pdf/2410.02703",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "https://arxiv.org/abs/2410.05258",
    2) | Out-Null

# --- Remove the four paragraphs that sat between the skepticism paragraph
#     and the (now relocated-in-place) link paragraph: "אבל מה זה מטריצת F...",
#     "השיטה המוצעת...", "בנוסף כבר...", "יש לי תחושה..." ---
$start = $d.Paragraphs.Item(6).Range.Start
$end = $d.Paragraphs.Item(9).Range.End
$r = $d.Range($start, $end)
$r.Delete()
